# Fix test problem on importing xlsx: the given/family name values in
# row 2 (A2/B2) were swapped relative to the header row, so swap them
# back into the correct columns.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$valA2 = $ws.Range("A2").Value2
$valB2 = $ws.Range("B2").Value2

$ws.Range("A2").Value2 = $valB2
$ws.Range("B2").Value2 = $valA2

# Update the active selection to match the edited workbook state.
$ws.Range("B6").Select()
